# Add new "Pais"/"Estado" and "Cidade"/"Bairro" lookup tables below the
# existing data on Planilha1 (rows 105-131).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- formatting helpers ---------------------------------------------------
# C105 mirrors the lone underline-styled blank cell used elsewhere (e.g. C1).
$ws.Range("C1").Copy()
$ws.Range("C105").PasteSpecial(-4122)

# A106/C106 and A119/D119 reuse the shaded "header" look already used for
# the ID_Supermercado/ID_produto table header in row 93.
$ws.Range("A93").Copy()
$ws.Range("A106").PasteSpecial(-4122)
$ws.Range("C106").PasteSpecial(-4122)
$ws.Range("A119").PasteSpecial(-4122)
$ws.Range("D119").PasteSpecial(-4122)

# --- Pais column (A106:A117) ----------------------------------------------
$ws.Range("A106").Value = 'Pais'
$ws.Range("A107").Value = 'Brasil'
$ws.Range("A108").Value = 'EUA'
$ws.Range("A109").Value = 'França'
$ws.Range("A110").Value = 'Itália'
$ws.Range("A111").Value = 'Canadá'
$ws.Range("A112").Value = 'Chile'
$ws.Range("A113").Value = 'Russia'
$ws.Range("A114").Value = 'Alemanha'
$ws.Range("A115").Value = 'Mexico'
$ws.Range("A116").Value = 'gito'
$ws.Range("A117").Value = 'China'

# --- Estado column (C106:C116) --------------------------------------------
$ws.Range("C106").Value = 'Estado'
$ws.Range("C107").Value = 'Epirito santo'
$ws.Range("C108").Value = 'São Paulo'
$ws.Range("C109").Value = 'Goias'
$ws.Range("C110").Value = 'Amazonas'
$ws.Range("C111").Value = 'Mato grosso do sul'
$ws.Range("C112").Value = 'Rio de janeiro'
$ws.Range("C113").Value = 'Texas'
$ws.Range("C114").Value = 'California'
$ws.Range("C115").Value = 'Flórida'
$ws.Range("C116").Value = 'Alasca'

# --- Cidade/Bairro headers --------------------------------------------------
$ws.Range("A119").Value = 'Cidade'
$ws.Range("D119").Value = 'Bairro'

# --- Cidade column (A120:A131) ---------------------------------------------
$ws.Range("A120").Value = 'Serra'
$ws.Range("A121").Value = 'Vitoria'
$ws.Range("A122").Value = 'Vila Velha'
$ws.Range("A123").Value = 'Cariacica'
$ws.Range("A124").Value = 'Guarapari'
$ws.Range("A125").Value = 'Rio de janeiro'
$ws.Range("A126").Value = 'Salvador'
$ws.Range("A127").Value = 'Feira de santana'
$ws.Range("A128").Value = 'São paulo'
$ws.Range("A129").Value = 'Campinas'
$ws.Range("A130").Value = 'Rio de janeiro'
$ws.Range("A131").Value = 'Niterói'

# --- Bairro column (D120:D131) ----------------------------------------------
$ws.Range("D120").Value = 'Manguinhos'
$ws.Range("D121").Value = 'Laranjeiras'
$ws.Range("D122").Value = 'Goiabeiras'
$ws.Range("D123").Value = 'Jardim marilandia'
$ws.Range("D124").Value = 'Barceçpma'
$ws.Range("D125").Value = 'Feu rosa'
$ws.Range("D126").Value = 'São pedro'
$ws.Range("D127").Value = 'São torquato'
$ws.Range("D128").Value = 'Dombosco'
$ws.Range("D129").Value = 'Colina de laranjeiras'
$ws.Range("D130").Value = 'Marilândia'
$ws.Range("D131").Value = 'Coqueiral de itapuã'

# --- remaining Id/ID header cells -------------------------------------------
$ws.Range("D106").Value = 'Id'
$ws.Range("E106").Value = 'Id_pais'
$ws.Range("B119").Value = 'ID'
$ws.Range("C119").Value = 'id_ESTADO'
$ws.Range("E119").Value = 'id'
$ws.Range("F119").Value = 'id_cidade'

# --- sheet-level cosmetics ---------------------------------------------------
$ws.Columns("D").ColumnWidth = 18.7

# Scroll the view down so row 104 is at the top and select C114, matching
# where the author was working when the sheet was saved.
$ws.Range("C114").Select()
$excel.ActiveWindow.ScrollRow = 104
$excel.ActiveWindow.ScrollColumn = 1
